$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that used to hold only the lecturer's name ("5840671 - Francisco
# José Moreira Chaves" in B13/C13, no label in A13) is removed outright;
# everything below shifts up one row.
$ws.Rows("13:13").Delete()

$docente = "5840671 - Francisco José Moreira Chaves"

# Objetivos: value now just holds the lecturer identification instead of
# the long mission statement.
$ws.Range("B10").Value = $docente
$ws.Range("C10").Value = $docente

# Programa resumido: value becomes "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Programa: value becomes the activation date "01/01/2018".
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

# Método: value (previously empty after the row shift) now holds the
# lecturer identification again.
$ws.Range("B18").Value = $docente
$ws.Range("C18").Value = $docente

# Critério: value becomes the former "Método:" evaluation text.
$ws.Range("B19").Value = "As avaliações serão por meio de trabalhos em equipes ou provas individuais, conforme adequação ao conteúdo programático."
$ws.Range("C19").Value = "As avaliações serão por meio de trabalhos em equipes ou provas individuais, conforme adequação ao conteúdo programático."

# Norma de recuperação: value becomes the former "Critério:" text.
$ws.Range("B20").Value = "Serão aplicadas duas avaliações para compor a média que será a soma das duas provas, sendo o resultado dividido por dois."
$ws.Range("C20").Value = "Serão aplicadas duas avaliações para compor a média que será a soma das duas provas, sendo o resultado dividido por dois."

# Bibliografia: value becomes the former "Norma de recuperação:" text
# (the old bibliography list itself is gone, having been on the deleted
# trailing row).
$ws.Range("B21").Value = "A nota final será composta pela média obtida da nota do período somada à nota de recuperação"
$ws.Range("C21").Value = "A nota final será composta pela média obtida da nota do período somada à nota de recuperação"
